$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the descriptive text in C6 (sequence continues as powers of two)
$ws.Range("C6").Value = "The columns are set as [1, 2, 4, 8, 16, 32, 64, 128, 256, 512]"

# Narrow spacer column (J / column 10) width to 0.5 characters
$ws.Columns.Item(10).ColumnWidth = -0.33

# Move the active selection to A1
[void]$ws.Range("A1").Select()
